$d = $word.ActiveDocument

# Update the title/date line
$d.Content.Find.Execute("2025-07-21 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-22 Tuesday", 2)

# Update each table cell's divisor problem text, targeted by cell to avoid
# collisions between old/new values that overlap across cells.
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Find.Execute("77÷4=19, 1", $true, $false, $false, $false, $false, $true, 1, $false, "38÷2=19, 0", 2)
$tbl.Cell(1, 2).Range.Find.Execute("77÷7=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "46÷3=15, 1", 2)
$tbl.Cell(1, 3).Range.Find.Execute("60÷9=6, 6", $true, $false, $false, $false, $false, $true, 1, $false, "96÷4=24, 0", 2)
$tbl.Cell(1, 4).Range.Find.Execute("33÷3=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "22÷4=5, 2", 2)
$tbl.Cell(1, 5).Range.Find.Execute("50÷5=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "86÷5=17, 1", 2)
$tbl.Cell(5, 1).Range.Find.Execute("57÷6=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "24÷9=2, 6", 2)
$tbl.Cell(5, 2).Range.Find.Execute("13÷6=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "34÷2=17, 0", 2)
$tbl.Cell(5, 3).Range.Find.Execute("18÷6=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "73÷3=24, 1", 2)
$tbl.Cell(5, 4).Range.Find.Execute("16÷6=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "69÷7=9, 6", 2)
$tbl.Cell(5, 5).Range.Find.Execute("32÷2=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "82÷7=11, 5", 2)
$tbl.Cell(9, 1).Range.Find.Execute("57÷3=19, 0", $true, $false, $false, $false, $false, $true, 1, $false, "74÷3=24, 2", 2)
$tbl.Cell(9, 2).Range.Find.Execute("14÷2=7, 0", $true, $false, $false, $false, $false, $true, 1, $false, "28÷6=4, 4", 2)
$tbl.Cell(9, 3).Range.Find.Execute("98÷7=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "87÷8=10, 7", 2)
$tbl.Cell(9, 4).Range.Find.Execute("37÷4=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "10÷2=5, 0", 2)
$tbl.Cell(9, 5).Range.Find.Execute("65÷3=21, 2", $true, $false, $false, $false, $false, $true, 1, $false, "18÷6=3, 0", 2)
$tbl.Cell(13, 1).Range.Find.Execute("76÷9=8, 4", $true, $false, $false, $false, $false, $true, 1, $false, "36÷4=9, 0", 2)
$tbl.Cell(13, 2).Range.Find.Execute("61÷2=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "38÷2=19, 0", 2)
$tbl.Cell(13, 3).Range.Find.Execute("79÷6=13, 1", $true, $false, $false, $false, $false, $true, 1, $false, "14÷4=3, 2", 2)
$tbl.Cell(13, 4).Range.Find.Execute("39÷9=4, 3", $true, $false, $false, $false, $false, $true, 1, $false, "70÷3=23, 1", 2)
$tbl.Cell(13, 5).Range.Find.Execute("67÷6=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "35÷8=4, 3", 2)
$tbl.Cell(17, 1).Range.Find.Execute("44÷7=6, 2", $true, $false, $false, $false, $false, $true, 1, $false, "37÷6=6, 1", 2)
$tbl.Cell(17, 2).Range.Find.Execute("19÷9=2, 1", $true, $false, $false, $false, $false, $true, 1, $false, "38÷2=19, 0", 2)
$tbl.Cell(17, 3).Range.Find.Execute("67÷2=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "49÷5=9, 4", 2)
$tbl.Cell(17, 4).Range.Find.Execute("42÷9=4, 6", $true, $false, $false, $false, $false, $true, 1, $false, "51÷4=12, 3", 2)
$tbl.Cell(17, 5).Range.Find.Execute("22÷8=2, 6", $true, $false, $false, $false, $false, $true, 1, $false, "52÷3=17, 1", 2)
